{"js": "// Replace each target text run with its updated value.\n// Every old value in this document is unique, so a body.search()\n// per pair unambiguously finds the single run to update.\nconst replacements = [\n  [\"2024-12-24 Tuesday\", \"2024-12-25 Wednesday\"],\n  [\"68\u00d730=\", \"35\u00d740=\"],\n  [\"12\u00d745=\", \"43\u00d799=\"],\n  [\"88\u00d726=\", \"49\u00d795=\"],\n  [\"53\u00d767=\", \"82\u00d774=\"],\n  [\"65\u00d793=\", \"94\u00d726=\"],\n  [\"63\u00d732=\", \"54\u00d735=\"],\n  [\"62\u00d718=\", \"43\u00d776=\"],\n  [\"22\u00d786=\", \"87\u00d757=\"],\n  [\"41\u00d775=\", \"13\u00d791=\"],\n  [\"76\u00d772=\", \"92\u00d790=\"],\n  [\"53\u00d725=\", \"99\u00d789=\"],\n  [\"71\u00d785=\", \"80\u00d742=\"],\n  [\"67\u00d758=\", \"38\u00d741=\"],\n  [\"33\u00d755=\", \"15\u00d726=\"],\n  [\"92\u00d781=\", \"26\u00d786=\"],\n  [\"36\u00d763=\", \"44\u00d726=\"],\n  [\"96\u00d767=\", \"55\u00d718=\"],\n  [\"40\u00d738=\", \"94\u00d749=\"],\n  [\"12\u00d722=\", \"61\u00d777=\"],\n  [\"74\u00d741=\", \"60\u00d794=\"],\n  [\"89\u00d719=\", \"67\u00d794=\"],\n  [\"47\u00d755=\", \"80\u00d741=\"],\n  [\"44\u00d725=\", \"68\u00d773=\"],\n  [\"90\u00d794=\", \"31\u00d777=\"],\n  [\"21\u00d747=\", \"96\u00d752=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each target text run with its updated value via Find/Replace.\n# Every old value in this document is unique, so FindText uniquely\n# identifies the single run to update for each pair.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-12-24 Tuesday\", \"2024-12-25 Wednesday\"),\n    @(\"68\u00d730=\", \"35\u00d740=\"),\n    @(\"12\u00d745=\", \"43\u00d799=\"),\n    @(\"88\u00d726=\", \"49\u00d795=\"),\n    @(\"53\u00d767=\", \"82\u00d774=\"),\n    @(\"65\u00d793=\", \"94\u00d726=\"),\n    @(\"63\u00d732=\", \"54\u00d735=\"),\n    @(\"62\u00d718=\", \"43\u00d776=\"),\n    @(\"22\u00d786=\", \"87\u00d757=\"),\n    @(\"41\u00d775=\", \"13\u00d791=\"),\n    @(\"76\u00d772=\", \"92\u00d790=\"),\n    @(\"53\u00d725=\", \"99\u00d789=\"),\n    @(\"71\u00d785=\", \"80\u00d742=\"),\n    @(\"67\u00d758=\", \"38\u00d741=\"),\n    @(\"33\u00d755=\", \"15\u00d726=\"),\n    @(\"92\u00d781=\", \"26\u00d786=\"),\n    @(\"36\u00d763=\", \"44\u00d726=\"),\n    @(\"96\u00d767=\", \"55\u00d718=\"),\n    @(\"40\u00d738=\", \"94\u00d749=\"),\n    @(\"12\u00d722=\", \"61\u00d777=\"),\n    @(\"74\u00d741=\", \"60\u00d794=\"),\n    @(\"89\u00d719=\", \"67\u00d794=\"),\n    @(\"47\u00d755=\", \"80\u00d741=\"),\n    @(\"44\u00d725=\", \"68\u00d773=\"),\n    @(\"90\u00d794=\", \"31\u00d777=\"),\n    @(\"21\u00d747=\", \"96\u00d752=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        $find.Text,        # FindText\n        $false,            # MatchCase\n        $false,            # MatchWholeWord\n        $false,            # MatchWildcards\n        $false,            # MatchSoundsLike\n        $false,            # MatchAllWordForms\n        $true,             # Forward\n        1,                 # Wrap (wdFindContinue)\n        $false,            # Format\n        $find.Replacement.Text,  # ReplaceWith\n        2                  # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n"}
